# Fix Training Data Issue (#48)
# The "Date" column (BF) held values in the wrong shape ("4-23-2007-08")
# because of how NBA stats were originally exported; correct them to the
# proper ISO date-as-text value "2008-04-23" for every data row (rows 2-31).
#
# NumberFormatLocal is set to the Text format ("@") before the value is
# written so the ISO-looking string "2008-04-23" is stored as literal text
# instead of being auto-parsed into a date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 31
$col = "BF"

$dateRange = $ws.Range("$col$firstRow`:$col$lastRow")
$dateRange.NumberFormatLocal = "@"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Range("$col$row").Value = "2008-04-23"
}
